$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D9").Value = "대통령 선거로 만든 Data Science 문제"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/presidential-election-problem-set/#utm_source=rss&utm_medium=rss&utm_campaign=presidential-election-problem-set"

$ws.Range("D26").Value = "생성 모델의 새로운 흐름 확산 모델(Diffusion model)에 관하여"

$ws.Range("D37").Value = "[Paper Review] Data diversification: A Simple Strategy For Neural Machine Translation"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1970&mod=document&pageid=1"
